# Commit: "update practical distributions and sensitivity with 5 SD for LI"
#
# Adds a new "iqr" (inter-quartile range) row as the last row of the
# descriptive-statistics table on every worksheet of the workbook.
# Before the edit each sheet has data in rows 1-8 (header + 7 stat rows,
# n_samples/min/max/mean/median/q25/q75). This appends row 9 with the
# "iqr" label in column A and the corresponding numeric values in B/C/D.

$wb = $excel.ActiveWorkbook

$rowsToAdd = @(
    @{ Sheet = "Maissilage";          B = 44.770931233286518;  C = 6.030423844792228;   D = 6.914877646550174 },
    @{ Sheet = "Grassilage";          B = 66.161040167890917;  C = 4.3830038409526111;  D = 11.88238002467928 },
    @{ Sheet = "Getreidestroh";       B = 21.939167145063806;  C = 2.1966571813689848;  D = 2.4068238448213286 },
    @{ Sheet = "Zuckerruebensilage";  B = 325.46881114621226;  C = 4.2158668272912818;  D = 0.17859390688718457 },
    @{ Sheet = "Rinderguelle";        B = 6.2344535282770224;  C = 1.5592889140136439;  D = 3.15953085543519 },
    @{ Sheet = "Schweineguelle";      B = 8.985642530511484;   C = 11.357000902078557;  D = 3.4844534329008301 },
    @{ Sheet = "HTK";                 B = 78.852603409405461;  C = 75.041473619582192;  D = 9.9083796591708051 }
)

foreach ($entry in $rowsToAdd) {
    $ws = $wb.Worksheets.Item($entry.Sheet)

    $ws.Cells.Item(9, 1).Value2 = "iqr"
    $ws.Cells.Item(9, 2).Value2 = $entry.B
    $ws.Cells.Item(9, 3).Value2 = $entry.C
    $ws.Cells.Item(9, 4).Value2 = $entry.D
}
